# Update the "取得日時" (acquired timestamp) column for all data rows
# on the "ランサーズ" sheet from 2025-12-21 01:38:23 to 2025-12-21 02:02:58.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-21 02:02:58"

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
